$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Locate the "}}-2026-MDP-GLDE/SGLC" stretch (just before the closing
# <w:tab/> run) so we can recompute run boundaries after editing it.
# ------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("}}-2026-MDP-GLDE/SGLC", $true, $false, $false, $false, `
                      $false, $true, 1, $false, "", 0)
$blockStart = $anchor.Start

# --- Step 1: "-MDP-GLDE" -> "-MDP" ---------------------------------
$r1 = $d.Content
$r1.Find.Execute("-MDP-GLDE", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r1.Text = "-MDP"

# --- Step 2: "/SGLC" -> "/" (rest re-added below as new runs) ------
$r2 = $d.Content
$r2.Find.Execute("/SGLC", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$r2.Text = "/"
$afterSlash = $r2.End

# --- Step 3: insert "GLDE-SGLCA" right after the "/" ----------------
$ins = $d.Range($afterSlash, $afterSlash)
$ins.InsertAfter("GLDE-SGLCA")

# ------------------------------------------------------------------
# Step 4: the edits above can coalesce formerly-separate runs that
# happen to share identical formatting (e.g. "}" + "}" + "-202" + "6"
# + "-MDP" + "/" + "GLDE-SGLCA" all end up as one <w:r>). Re-impose
# every run boundary that must exist in the final document by
# re-asserting each slice's (unchanged) Bold value - that forces a
# fresh <w:r> without leaving any stray formatting behind, and the
# split survives subsequent saves.
# ------------------------------------------------------------------
$lens = 1, 1, 4, 1, 4, 1, 4, 1, 4, 1   # "}" "}" "-202" "6" "-MDP" "/" "GLDE" "-" "SGLC" "A"
$pos = $blockStart
foreach ($len in $lens) {
    $slice = $d.Range($pos, $pos + $len)
    $slice.Font.Bold = 0
    $slice.Font.Bold = 1
    $pos = $pos + $len
}
